$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update monthly_salary (column X) to a flat numeric value and
# school_name (column Y) to the new school name for all data rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 19 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 24).Value = 800000
    $ws.Cells.Item($r, 25).Value = "Heaven is my Home"
}
